{"js": "// Word Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// The target document is a fixed-width \"billing slip\" where every line is a\n// single paragraph/run. The edit swaps out retailer/store/amount details for\n// a new transaction while keeping every field's column width (and therefore\n// the surrounding run/paragraph formatting) identical.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// Map of exact current paragraph text -> exact replacement paragraph text.\n// (All replacement strings keep the same overall length as the originals,\n// matching the fixed-width layout used throughout the document.)\nconst replacements = new Map([\n  [\n    \"  ARIYAMANGALAM                                      Retailer Name   :  ESHWAR MEDICALS                    \",\n    \"  ARIYAMANGALAM                                      Retailer Name   :  SUBA MALIGAI-D                     \",\n  ],\n  [\n    \"  TRICHY-620010                                      Address         :  Thamarai 1st St Ezhil Nagar Tamil  \",\n    \"  TRICHY-620010                                      Address         :  GANDHI SALAI                       \",\n  ],\n  [\n    \"  PHONE NO         :9944951444                                          India Near Kumutha Store           \",\n    \"  PHONE NO         :9944951444                                          SOUTH KATTUR                       \",\n  ],\n  [\n    \"  GSTIN No         :33AAPFD1365C1ZR                                                                        \",\n    \"  GSTIN No         :33AAPFD1365C1ZR                                     TRICHY                             \",\n  ],\n  [\n    \"  RS PAN No        :AAPFD1365C                       Phone No        :   8072005857                        \",\n    \"  RS PAN No        :AAPFD1365C                       Phone No        :   8056384501                        \",\n  ],\n  [\n    \"  Salesperson Name :SAKTHIVEL M                      \",\n    \"  Salesperson Name :SAGARIAVINCENT                   \",\n  ],\n  [\n    \"  Beat Name        :Chemist - Thiruvarambur NUTS     GSTIN NO        :                                     \",\n    \"  Beat Name        :D-KATTUR-1  3S                   GSTIN NO        :                                     \",\n  ],\n  [\n    \"  HUL STORE ID     :HUL-41A392D-P25120               Time of Billing :   14/12/2023 22:26:36               \",\n    \"  HUL STORE ID     :HUL-414006D-P5185                Time of Billing :   26/12/2023 16:01:20               \",\n  ],\n  [\n    \"  Six Hundred Twenty-Nine Rupees Only                                        \",\n    \"  Two Thousand Six Hundred Forty-Six Rupees Only                             \",\n  ],\n  [\n    \"  ABC54009     ESHWAR MEDICALS    Amt : 629.00\",\n    \"  ABC56336     SUBA MALIGAI-D    Amt : 2646.00\",\n  ],\n]);\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const newText = replacements.get(para.text);\n  if (newText !== undefined) {\n    // Replace in place so paragraph/run formatting (e.g. the bold, sz=24\n    // run on the \"ABC...\" summary line) is preserved.\n    para.getRange().insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# The target document is a fixed-width \"billing slip\" where every line is a\n# single paragraph/run. The edit swaps out retailer/store/amount details for\n# a new transaction while keeping every field's column width (and therefore\n# the surrounding run/paragraph formatting) identical.\n\n$d = $word.ActiveDocument\n\n# Map of exact current paragraph text (paragraph mark excluded) -> exact\n# replacement paragraph text. (All replacement strings keep the same overall\n# length as the originals, matching the fixed-width layout used throughout\n# the document.)\n$map = @{\n  \"  ARIYAMANGALAM                                      Retailer Name   :  ESHWAR MEDICALS                    \" = \"  ARIYAMANGALAM                                      Retailer Name   :  SUBA MALIGAI-D                     \"\n  \"  TRICHY-620010                                      Address         :  Thamarai 1st St Ezhil Nagar Tamil  \" = \"  TRICHY-620010                                      Address         :  GANDHI SALAI                       \"\n  \"  PHONE NO         :9944951444                                          India Near Kumutha Store           \" = \"  PHONE NO         :9944951444                                          SOUTH KATTUR                       \"\n  \"  GSTIN No         :33AAPFD1365C1ZR                                                                        \" = \"  GSTIN No         :33AAPFD1365C1ZR                                     TRICHY                             \"\n  \"  RS PAN No        :AAPFD1365C                       Phone No        :   8072005857                        \" = \"  RS PAN No        :AAPFD1365C                       Phone No        :   8056384501                        \"\n  \"  Salesperson Name :SAKTHIVEL M                      \"                                                       = \"  Salesperson Name :SAGARIAVINCENT                   \"\n  \"  Beat Name        :Chemist - Thiruvarambur NUTS     GSTIN NO        :                                     \" = \"  Beat Name        :D-KATTUR-1  3S                   GSTIN NO        :                                     \"\n  \"  HUL STORE ID     :HUL-41A392D-P25120               Time of Billing :   14/12/2023 22:26:36               \" = \"  HUL STORE ID     :HUL-414006D-P5185                Time of Billing :   26/12/2023 16:01:20               \"\n  \"  Six Hundred Twenty-Nine Rupees Only                                        \"                               = \"  Two Thousand Six Hundred Forty-Six Rupees Only                             \"\n  \"  ABC54009     ESHWAR MEDICALS    Amt : 629.00\"                                                               = \"  ABC56336     SUBA MALIGAI-D    Amt : 2646.00\"\n}\n\nforeach ($p in $d.Paragraphs) {\n    # Paragraph.Range.Text includes the trailing paragraph mark (and, for the\n    # very last paragraph, a section mark) -- strip those control chars\n    # before comparing so the lookup matches the diff's paragraph text.\n    $full = $p.Range.Text\n    $key = $full.TrimEnd([char]13, [char]7)\n    if ($map.ContainsKey($key)) {\n        $r = $p.Range\n        # Only touch the text portion, leaving the paragraph mark alone, so\n        # paragraph/run formatting (bold, font size, etc.) is preserved.\n        $r.SetRange($r.Start, $r.Start + $key.Length)\n        $r.Text = $map[$key]\n    }\n}\n"}
